$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

function Find-ParagraphIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# 1) Title paragraph: "Inicidencias Honcizek" -> "Inicidencias" + " " + "Honcizek",
#    each word wrapped in proofErr spellStart/spellEnd (the space run is left unwrapped).
$titleIdx = Find-ParagraphIndex $d "Inicidencias Honcizek"
$pTitle = $d.Paragraphs.Item($titleIdx)
$titleXml = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="40"/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="40"/></w:rPr><w:t>Inicidencias</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="40"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="40"/></w:rPr><w:t>Honcizek</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
$pTitle.Range.InsertXML($titleXml)

# 2) "Admin" paragraph: wrap the existing run in proofErr spellStart/spellEnd.
$adminIdx = Find-ParagraphIndex $d "Admin"
$pAdmin = $d.Paragraphs.Item($adminIdx)
$adminXml = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:pPr><w:rPr><w:b/><w:sz w:val="40"/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="40"/></w:rPr><w:t>Admin</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
$pAdmin.Range.InsertXML($adminXml)

# 3) Remove the "Suscripciones" heading paragraph entirely.
$suscripcionesIdx = Find-ParagraphIndex $d "Suscripciones"
$d.Paragraphs.Item($suscripcionesIdx).Range.Delete()

# 4) Collapse the bullet paragraph ("Rellenar select de proyectos...") together with the
#    empty paragraph that immediately follows it into a single plain paragraph that keeps
#    only the _GoBack bookmark (the bullet text and its list formatting are dropped).
$bulletIdx = Find-ParagraphIndex $d "Rellenar select de proyectos"
$pBullet = $d.Paragraphs.Item($bulletIdx)
$pTrailing = $d.Paragraphs.Item($bulletIdx + 1)
$mergeRange = $d.Range($pBullet.Range.Start, $pTrailing.Range.End)
$mergeXml = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
$mergeRange.InsertXML($mergeXml)
